# Update column F (dSF) values for rows 4-27 per repulled/recomputed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 9
    5  = 4
    6  = 4
    7  = 4
    8  = 1
    9  = 1
    10 = 2
    11 = -1
    12 = 2
    13 = 9
    14 = 2
    15 = 1
    16 = -4
    17 = 6
    19 = -2
    20 = 5
    22 = 7
    23 = 4
    24 = 3
    26 = 2
    27 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
